$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift header row: Resto, Menu, Quantity, Harga
$ws.Range("A1").Value = "Resto"
$ws.Range("B1").Value = "Menu"
$ws.Range("C1").Value = "Quantity"
$ws.Range("D1").Value = "Harga"

# Data row
$ws.Range("A2").Value = "Oto Bento"
$ws.Range("B2").Value = "Chicken Blackpepper"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 25000

# Column widths
$ws.Columns.Item(4).ColumnWidth = 20.42578125
$ws.Columns.Item(5).ColumnWidth = 13.85546875

# Selection
$ws.Range("E6").Select()
